# The title block in the document currently reads "...110317" as a
# single run (date stamp 11/03/17). The update changes the date stamp to
# "050517" (05/05/17) and, per the target OOXML, splits it across two
# runs: "0505" (keeping the original run) followed by a new run "17"
# that carries the same run formatting (Times New Roman, bold, kern 36,
# sz/szCs 16).

$d = $word.ActiveDocument

# Locate "110317" and remember its extent.
$hit = $d.Content
$found = $hit.Find.Execute("110317", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find '110317' in the document"
}
$start = $hit.Start
$end = $hit.End

# First four characters ("1103") become "0505"; keep this inside the
# range that used to hold the whole date stamp, so it stays the first
# (original) run.
$rngFirst = $d.Range($start, $start + 4)
$rngFirst.Text = "0505"

# Remaining two characters ("17") are re-typed into their own range,
# which Word materialises as a brand-new run right after the first one.
$rngSecond = $d.Range($start + 4, $start + 6)
$rngSecond.Text = "17"

# The two runs now share identical formatting and are adjacent, so make
# sure the formatting actually matches (it already does by inheritance,
# this just makes the intent explicit and keeps them from accidentally
# differing), without disturbing the just-created run boundary.
$origSize = $rngFirst.Font.Size
$rngFirst.Font.Size = $origSize + 1
$rngFirst.Font.Size = $origSize
